$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update).
# D-column numeric-looking price strings must stay text (not be
# coerced to numbers) to match the source data, which stores them
# as plain display strings (e.g. "1.00", "2.40" must keep trailing zeros).

# Row 2
$ws.Cells.Item(2, 4).Value = '42.153.70'
$ws.Cells.Item(2, 5).Value = '  +5.13%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.243.60'
$ws.Cells.Item(3, 5).Value = '  +1.72%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '232.49'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.61%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.23%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '61.65'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -3.70%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.10%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +1.76%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '59.41'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.69%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0899'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +4.42%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.18%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '2.576.66'
$ws.Cells.Item(13, 5).Value = '  +1.55%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '15.76'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -2.33%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '22.18'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.21%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.807'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.59%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.10%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.265.60'
$ws.Cells.Item(18, 5).Value = '  +2.42%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '42.159.08'
$ws.Cells.Item(19, 5).Value = '  +5.27%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0906'
$ws.Cells.Item(20, 5).Value = '  -0.65%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '72.43'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.08%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.99%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '253.08'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +8.98%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.07%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'PancakeSwap'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.40'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.86%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.32'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.35%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.69'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.02%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.78%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '168.92'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.74%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '20.09'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.21%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -3.33%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.72'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.32%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.67%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +5.88%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.67'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.47%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0639'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.24%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.69'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -5.44%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.73'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -3.94%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -4.00%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.000258'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +31.67%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.15%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +4.74%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -2.55%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +3.98%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.26%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0975'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +5.03%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '99.39'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.28%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '1.482.20'
$ws.Cells.Item(48, 5).Value = '  -2.67%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '16.57'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -7.17%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.20%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -2.49%  '
